$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new row of data (row 54): Rotate Array ---
$ws.Range("A54").Value = "Rotate Array"
$ws.Range("B54").Value = "Array"
$ws.Range("C54").Value = "No"
$ws.Range("D54").Value = "Yes"
$ws.Range("E54").Value = "Medium"
$ws.Range("F54").Value = "Medium"

# Hyperlink + display text for G54 (matches style used by other problem links)
$ws.Hyperlinks.Add($ws.Range("G54"), "189 - Rotate Array", "", "", "189 - Rotate Array")
$ws.Range("G54").Style = $ws.Range("G53").Style

# --- Update the view: selection moved to L7, no forced top-left scroll ---
$ws.Range("L7").Select()

# --- Extend data validation ranges down to row 54 (recreate in original order) ---
$ws.Range("E2:F53").Validation.Delete()
$ws.Range("C2:C53").Validation.Delete()
$ws.Range("D2:D53").Validation.Delete()
$ws.Range("B2:B53").Validation.Delete()

$ws.Range("E2:F54").Validation.Add(3, 1, 1, '"Easy, Medium, Hard"')
$ws.Range("C2:C54").Validation.Add(3, 1, 1, '"Yes, No"')
$ws.Range("C2:C54").Validation.IgnoreBlank = $false
$ws.Range("D2:D54").Validation.Add(3, 1, 1, '"Yes, No"')
$ws.Range("B2:B54").Validation.Add(3, 1, 1, '"Array, Binary, Dynamic Programming, Graph, Interval, Linked List, Matrix, String, Tree, Heap, Class Design"')
